$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: turn the former last row into a "continuation" row (like rows 4/7) ---
# Copy the formatting (styles + bottom border) from row 4 onto row 9, keeping
# row 9's own A9 value (SCRIPT/T01P01A/um1113.ssb) intact; B9:E9 stay empty.
$ws.Range("A4:E4").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)

# --- Row 10: new data row (SCRIPT/T01P02A/um1319.ssb, line 83) ---
# Copy formatting from row 3 (a normal, non-continuation data row).
$ws.Range("A3:E3").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

# --- Row 11: new data row (SCRIPT/T01P02A/um1405.ssb, line 86) ---
$ws.Range("A3:E3").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

# Fill in the new cell values in the same order the strings were first
# entered (keeps shared-string table order stable / matches source order).
$ws.Range("A10").Value = "SCRIPT/T01P02A/um1319.ssb"
$ws.Range("C10").Value = ' Time stopping...[K]might not make\nany difference to us.'
$ws.Range("C11").Value = ' We\''re not moving either way.'
$ws.Range("D10").Value = ' Остановка времени...[K] Для нас\nникакой разницы не будет.'
$ws.Range("D11").Value = ' Мы и так не двигаемся.'
$ws.Range("E10").Value = ' Ïòóàîïâëà âñåíåîé...[K] Äìÿ îàò\nîéëàëïê ñàèîéøú îå áôäåó.'
$ws.Range("E11").Value = ' Íú é óàë îå äâéãàåíòÿ.'
$ws.Range("A11").Value = "SCRIPT/T01P02A/um1405.ssb"

$ws.Range("B10").Value = 83
$ws.Range("B11").Value = 86

# Rows 10/11 hold wrapped two-line strings, same row height as the other
# multi-line rows (e.g. row 9).
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Rows.Item(11).RowHeight = 43.2

# Update the active selection to D9, matching the edited workbook.
$ws.Range("D9").Select() | Out-Null
